$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8
$ws.Range("A8").Value = "[180, 147.27, 114.54, 81.81, 49.09, 16.36, 0, -16.36, -49.09, -81.81, -114.54, -147.27]"
$ws.Range("B8").Value = "[2]"
$ws.Range("C8").Value = "[0.04]"
$ws.Range("D8").Value = 5
$ws.Range("E8").Value = 10
$ws.Range("F8").Value = 6
$ws.Range("H8").Value = "long free"

# Row 9
$ws.Range("A9").Value = "[180, 147.27, 114.54, 81.81, 49.09, 16.36, 0, -16.36, -49.09, -81.81, -114.54, -147.27]"
$ws.Range("B9").Value = "[2]"
$ws.Range("C9").Value = "[0.04]"
$ws.Range("D9").Value = 5
$ws.Range("E9").Value = 1
$ws.Range("F9").Value = 6
$ws.Range("H9").Value = "shadow pause free"

# Update the selected cell in the sheet view
$ws.Range("F21").Select()
